$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.845.47"
$ws.Range("E2").Value = "  -0.88%  "
$ws.Range("D3").Value = "2.238.26"
$ws.Range("E3").Value = "  -0.29%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.63"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +7.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.633"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.21%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "72.11"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.95%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.566"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.97%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.24"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +18.52%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0974"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.46%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "58.39"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.41%  "
$ws.Range("E13").Value = "  +0.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.91"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.87%  "
$ws.Range("D15").Value = "2.572.98"
$ws.Range("E15").Value = "  -0.02%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.11"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.54%  "
$ws.Range("E17").Value = "  -1.60%  "
$ws.Range("D18").Value = "2.237.60"
$ws.Range("E18").Value = "  -1.04%  "
$ws.Range("D19").Value = "41.809.51"
$ws.Range("E19").Value = "  -0.64%  "
$ws.Range("E20").Value = "  -2.28%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.25"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.55%  "
$ws.Range("B22").Value = "Litecoin"
$ws.Range("C22").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.43"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.71"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.25%  "
$ws.Range("B24").Value = "ImmutableX"
$ws.Range("C24").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.25"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +25.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.75"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +3.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.51"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +5.85%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.15"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.20%  "
$ws.Range("E29").Value = "  +5.26%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "172.22"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +3.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.80"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.123"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.55%  "
$ws.Range("E33").Value = "  -1.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.48"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +3.63%  "
$ws.Range("E35").Value = "  +1.61%  "
$ws.Range("B36").Value = "Filecoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.72"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.99%  "
$ws.Range("B37").Value = "InjectiveProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "26.45"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +27.27%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.17"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +16.49%  "
$ws.Range("E39").Value = "  +5.59%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.29"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.82%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.06"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "68.55"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +6.50%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.216"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +15.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.04"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.38%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "11.69"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +22.17%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.88"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +12.84%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.88"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.32%  "
$ws.Range("E48").Value = "  +2.02%  "
$ws.Range("E49").Value = "  +0.12%  "
$ws.Range("D50").Value = "0.0₃0152"
$ws.Range("E50").Value = "  +19.47%  "
$ws.Range("E51").Value = "  +1.69%  "
